# Auto-generated edit script: refresh cached market-price / profit
# columns (H:N) across all 8 profession sheets, per the scheduled-runner
# data refresh described in the commit.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2793
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 3689.5
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 3689.5
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -3915.5
$ws.Range("H12").Value = 12772.125
$ws.Range("I12").Value = 14526
$ws.Range("J12").Value = 495
$ws.Range("K12").Value = 14526
$ws.Range("L12").Value = 495
$ws.Range("M12").Value = -14356
$ws.Range("N12").Value = -835
$ws.Range("H70").Value = 22125
$ws.Range("J70").Value = 28166.666
$ws.Range("L70").Value = 84499.99800000001
$ws.Range("N70").Value = -85039.99800000001
$ws.Range("H73").Value = 22125
$ws.Range("J73").Value = 28166.666
$ws.Range("L73").Value = 84499.99800000001
$ws.Range("N73").Value = -86371.99800000001
$ws.Range("H132").Value = 3166.0908
$ws.Range("I132").Value = 3215.2222
$ws.Range("K132").Value = 9645.6666
$ws.Range("M132").Value = -7115.6666

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4958.5586
$ws.Range("I32").Value = 2278.0527
$ws.Range("J32").Value = 12598
$ws.Range("K32").Value = 2278.0527
$ws.Range("L32").Value = 12598
$ws.Range("M32").Value = -1991.0527
$ws.Range("N32").Value = -13172
$ws.Range("H63").Value = 2854.4285
$ws.Range("I63").Value = 2476.4
$ws.Range("K63").Value = 2476.4
$ws.Range("M63").Value = -1790.4
$ws.Range("H66").Value = 2854.4285
$ws.Range("I66").Value = 2476.4
$ws.Range("K66").Value = 12382
$ws.Range("M66").Value = -8950
$ws.Range("H86").Value = 49999
$ws.Range("J86").Value = 49999
$ws.Range("L86").Value = 49999
$ws.Range("N86").Value = -52371
$ws.Range("H89").Value = 49999
$ws.Range("J89").Value = 49999
$ws.Range("L89").Value = 149997
$ws.Range("N89").Value = -161853
$ws.Range("H122").Value = 1918.2727
$ws.Range("J122").Value = 1623.875
$ws.Range("L122").Value = 4871.625
$ws.Range("N122").Value = -9771.625

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 492
$ws.Range("I94").Value = 459.5
$ws.Range("K94").Value = 459.5
$ws.Range("M94").Value = -8.5
$ws.Range("H134").Value = 22731676
$ws.Range("I134").Value = 29415112
$ws.Range("K134").Value = 88245336
$ws.Range("M134").Value = -88242801

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8154.7236
$ws.Range("J31").Value = 14711
$ws.Range("L31").Value = 14711
$ws.Range("N31").Value = -15301
$ws.Range("H34").Value = 8154.7236
$ws.Range("J34").Value = 14711
$ws.Range("L34").Value = 14711
$ws.Range("N34").Value = -15115
$ws.Range("H58").Value = 50024380
$ws.Range("I58").Value = 71458110
$ws.Range("K58").Value = 71458110
$ws.Range("M58").Value = -71457907
$ws.Range("H107").Value = 700340.4399999999
$ws.Range("I107").Value = 988663.75
$ws.Range("K107").Value = 988663.75
$ws.Range("M107").Value = -986743.75
$ws.Range("H135").Value = 99995.5
$ws.Range("J135").Value = 99995.5
$ws.Range("L135").Value = 99995.5
$ws.Range("N135").Value = -110135.5
$ws.Range("H136").Value = 50024380
$ws.Range("I136").Value = 71458110
$ws.Range("K136").Value = 214374330
$ws.Range("M136").Value = -214371780

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 59599.293
$ws.Range("I5").Value = 84137.914
$ws.Range("J5").Value = 706.6
$ws.Range("K5").Value = 252413.742
$ws.Range("L5").Value = 2119.8
$ws.Range("M5").Value = -252301.742
$ws.Range("N5").Value = -2343.8
$ws.Range("H12").Value = 181.24
$ws.Range("I12").Value = 169.72728
$ws.Range("J12").Value = 190.28572
$ws.Range("K12").Value = 509.18184
$ws.Range("L12").Value = 570.85716
$ws.Range("M12").Value = -336.18184
$ws.Range("N12").Value = -916.85716
$ws.Range("H121").Value = 949197.4
$ws.Range("I121").Value = 204195.4
$ws.Range("J121").Value = 1414823.6
$ws.Range("K121").Value = 612586.2
$ws.Range("L121").Value = 4244470.800000001
$ws.Range("M121").Value = -611276.2
$ws.Range("N121").Value = -4247090.800000001
$ws.Range("H135").Value = 59599.293
$ws.Range("I135").Value = 84137.914
$ws.Range("J135").Value = 706.6
$ws.Range("K135").Value = 757241.226
$ws.Range("L135").Value = 6359.400000000001
$ws.Range("M135").Value = -754706.226
$ws.Range("N135").Value = -11429.4
$ws.Range("H136").Value = 899
$ws.Range("I136").Value = 899
$ws.Range("K136").Value = 2697
$ws.Range("M136").Value = 2403

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2150.3333
$ws.Range("I3").Value = 966
$ws.Range("J3").Value = 3334.6667
$ws.Range("K3").Value = 966
$ws.Range("L3").Value = 3334.6667
$ws.Range("M3").Value = -850
$ws.Range("N3").Value = -3566.6667
$ws.Range("H7").Value = 1676000.5
$ws.Range("I7").Value = 5001000
$ws.Range("J7").Value = 13500.75
$ws.Range("K7").Value = 5001000
$ws.Range("L7").Value = 13500.75
$ws.Range("M7").Value = -5000888
$ws.Range("N7").Value = -13724.75
$ws.Range("H8").Value = 1676000.5
$ws.Range("I8").Value = 5001000
$ws.Range("J8").Value = 13500.75
$ws.Range("K8").Value = 5001000
$ws.Range("L8").Value = 13500.75
$ws.Range("M8").Value = -5000861
$ws.Range("N8").Value = -13778.75
$ws.Range("H12").Value = 6982.8335
$ws.Range("I12").Value = 4999
$ws.Range("J12").Value = 7379.6
$ws.Range("K12").Value = 4999
$ws.Range("L12").Value = 7379.6
$ws.Range("M12").Value = -4859
$ws.Range("N12").Value = -7659.6
$ws.Range("H14").Value = 1678666.9
$ws.Range("I14").Value = 2502249
$ws.Range("K14").Value = 2502249
$ws.Range("M14").Value = -2502081
$ws.Range("H97").Value = 2454.3333
$ws.Range("I97").Value = 2345.4
$ws.Range("J97").Value = 2999
$ws.Range("K97").Value = 2345.4
$ws.Range("L97").Value = 2999
$ws.Range("M97").Value = -1849.4
$ws.Range("N97").Value = -3991
$ws.Range("H107").Value = 1297.8334
$ws.Range("I107").Value = 958.75
$ws.Range("J107").Value = 1976
$ws.Range("K107").Value = 958.75
$ws.Range("L107").Value = 1976
$ws.Range("M107").Value = 961.25
$ws.Range("N107").Value = -5816
$ws.Range("H113").Value = 46916.086
$ws.Range("I113").Value = 57948.5
$ws.Range("K113").Value = 57948.5
$ws.Range("M113").Value = -55778.5
$ws.Range("H126").Value = 6363.9165
$ws.Range("I126").Value = 6733.364
$ws.Range("K126").Value = 20200.092
$ws.Range("M126").Value = -17730.092
$ws.Range("H132").Value = 6254821
$ws.Range("I132").Value = 7357053.5
$ws.Range("K132").Value = 22071160.5
$ws.Range("M132").Value = -22068630.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 13200
$ws.Range("J3").Value = 13200
$ws.Range("L3").Value = 13200
$ws.Range("N3").Value = -13424
$ws.Range("H14").Value = 28000
$ws.Range("J14").Value = 28000
$ws.Range("L14").Value = 28000
$ws.Range("N14").Value = -28344
$ws.Range("H15").Value = 13200
$ws.Range("J15").Value = 13200
$ws.Range("L15").Value = 13200
$ws.Range("N15").Value = -13540
$ws.Range("H61").Value = 2339.3901
$ws.Range("I61").Value = 2061.7715
$ws.Range("K61").Value = 2061.7715
$ws.Range("M61").Value = -1859.7715
$ws.Range("H113").Value = 2339.3901
$ws.Range("I113").Value = 2061.7715
$ws.Range("K113").Value = 2061.7715
$ws.Range("M113").Value = 108.2285000000002
$ws.Range("H132").Value = 53347664
$ws.Range("J132").Value = 3996.6667
$ws.Range("L132").Value = 11990.0001
$ws.Range("N132").Value = -17050.0001
$ws.Range("H136").Value = 2855.25
$ws.Range("I136").Value = 2567.5
$ws.Range("J136").Value = 2999.125
$ws.Range("K136").Value = 7702.5
$ws.Range("L136").Value = 8997.375
$ws.Range("M136").Value = -5152.5
$ws.Range("N136").Value = -14097.375

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H113").Value = 384.85
$ws.Range("J113").Value = 729.6
$ws.Range("L113").Value = 2188.8
$ws.Range("N113").Value = -6528.8
$ws.Range("H136").Value = 20834884
$ws.Range("I136").Value = 26316474
$ws.Range("J136").Value = 4836.4
$ws.Range("K136").Value = 78949422
$ws.Range("L136").Value = 14509.2
$ws.Range("M136").Value = -78946872
